$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header cells I1 and J1, matching the style of existing header row (style index 1)
$ws.Range("H1").Copy($ws.Range("I1:J1"))
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# Data for columns I (I0) and J (IF), rows 2-54
$data = @(
    ,@(2, 5, 6)
    ,@(3, 7, 8)
    ,@(4, 3, 4)
    ,@(5, 7, 7)
    ,@(6, 7, 7)
    ,@(7, 6, 6)
    ,@(8, 8, 8)
    ,@(9, 6, 6)
    ,@(10, 7, 7)
    ,@(11, 7, 7)
    ,@(12, 6, 6)
    ,@(13, 6, 6)
    ,@(14, 10, 10)
    ,@(15, 6, 6)
    ,@(16, 7, 7)
    ,@(17, 8, 8)
    ,@(18, 5, 6)
    ,@(19, 4, 6)
    ,@(20, 9, 9)
    ,@(21, 7, 7)
    ,@(22, 7, 7)
    ,@(23, 8, 8)
    ,@(24, 9, 9)
    ,@(25, 5, 7)
    ,@(26, 6, 7)
    ,@(27, 6, 7)
    ,@(28, 8, 8)
    ,@(29, 7, 8)
    ,@(30, 7, 7)
    ,@(31, 5, 5)
    ,@(32, 6, 6)
    ,@(33, 6, 7)
    ,@(34, 9, 9)
    ,@(35, 8, 8)
    ,@(36, 7, 7)
    ,@(37, 7, 7)
    ,@(38, 6, 6)
    ,@(39, 6, 6)
    ,@(40, 8, 8)
    ,@(41, 4, 5)
    ,@(42, 6, 7)
    ,@(43, 9, 9)
    ,@(44, 8, 8)
    ,@(45, 9, 9)
    ,@(46, 7, 8)
    ,@(47, 9, 9)
    ,@(48, 1, 1)
    ,@(49, 4, 5)
    ,@(50, 8, 8)
    ,@(51, 8, 8)
    ,@(52, 8, 8)
    ,@(53, 6, 6)
    ,@(54, 8, 8)
)

foreach ($item in $data) {
    $r = $item[0]
    $iVal = $item[1]
    $jVal = $item[2]
    $ws.Cells.Item($r, 9).Value2 = $iVal
    $ws.Cells.Item($r, 10).Value2 = $jVal
}

Write-Output "Applied I and J columns"
